$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above the current row 257 (shifts existing rows 257-309 down to 260-312)
$ws.Rows.Item(257).Insert()
$ws.Rows.Item(257).Insert()
$ws.Rows.Item(257).Insert()

# Common (unchanged) column values shared with the surrounding rows
$mercadoId = 11
$mercado = "Vega Monumental Concepción"
$region = "Bíobío"
$codreg = 8
$tipo = "Fruta"
$productoId = 100106
$producto = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"
$variedad = "Hass"

# Row 257: Hass, 1a nueva(o)
$r = 257
$ws.Cells.Item($r, 1).Value2 = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value2 = 44474
$ws.Cells.Item($r, 5).Value2 = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value2 = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value2 = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "1a nueva(o)"
$ws.Cells.Item($r, 13).Value2 = 50
$ws.Cells.Item($r, 14).Value2 = 2900
$ws.Cells.Item($r, 15).Value2 = 2900
$ws.Cells.Item($r, 16).Value2 = 2900
$ws.Cells.Item($r, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item($r, 18).Value = "Provincia de Quillota"
$ws.Cells.Item($r, 19).Value2 = 2900
$ws.Cells.Item($r, 20).Value2 = 1

# Row 258: Hass, 2a nueva(o)
$r = 258
$ws.Cells.Item($r, 1).Value2 = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value2 = 44474
$ws.Cells.Item($r, 5).Value2 = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value2 = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value2 = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "2a nueva(o)"
$ws.Cells.Item($r, 13).Value2 = 50
$ws.Cells.Item($r, 14).Value2 = 2700
$ws.Cells.Item($r, 15).Value2 = 2700
$ws.Cells.Item($r, 16).Value2 = 2700
$ws.Cells.Item($r, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item($r, 18).Value = "Provincia de Quillota"
$ws.Cells.Item($r, 19).Value2 = 2700
$ws.Cells.Item($r, 20).Value2 = 1

# Row 259: Hass, 3a nueva (o)
$r = 259
$ws.Cells.Item($r, 1).Value2 = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value2 = 44474
$ws.Cells.Item($r, 5).Value2 = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value2 = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value2 = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "3a nueva (o)"
$ws.Cells.Item($r, 13).Value2 = 50
$ws.Cells.Item($r, 14).Value2 = 2300
$ws.Cells.Item($r, 15).Value2 = 2300
$ws.Cells.Item($r, 16).Value2 = 2300
$ws.Cells.Item($r, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item($r, 18).Value = "Provincia de Quillota"
$ws.Cells.Item($r, 19).Value2 = 2300
$ws.Cells.Item($r, 20).Value2 = 1
